$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Axis Pattern")
$r = $ws.Range("C199")
$v = $r.QuotePrefix
Write-Host "QP=$v"
